$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Acknowledgments")

# 1. Delete the row for "lemmalist-greek" (alphabetically sorted table, row 10).
$targetRow = -1
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 1; $r -le $lastRow; $r++) {
    $val = $ws.Cells.Item($r, 1).Value()
    if ($val -eq "lemmalist-greek") {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    $ws.Rows.Item($targetRow).Delete()
}

# 2. Clear all existing hyperlinks on the sheet (deleting a row does not
#    renumber/remove the hyperlink collection in this runtime, so we rebuild
#    it from scratch to match the new row layout).
if ($ws.Hyperlinks.Count -gt 0) {
    $ws.Hyperlinks.Item(1).Range.Hyperlinks.Delete()
}

# 3. Re-create the hyperlinks with addresses shifted to account for the
#    deleted row, skipping the two hyperlinks that belonged to the removed row.
$hyperlinkData = @(
    @('B2', 'https://www.crummy.com/software/BeautifulSoup/', $null),
    @('B4', 'https://github.com/Ousret/charset_normalizer', $null),
    @('B8', 'https://github.com/Mimino666/langdetect', $null),
    @('B9', 'https://github.com/saffsd/langid.py', $null),
    @('B11', 'https://lxml.de/', $null),
    @('B16', 'https://www.numpy.org/', $null),
    @('B18', 'https://foss.heptapod.net/openpyxl/openpyxl', $null),
    @('B26', 'https://github.com/python-openxml/python-docx', $null),
    @('B28', 'https://github.com/psf/requests', $null),
    @('B30', 'https://scipy.org/scipylib/', $null),
    @('F2', 'https://bazaar.launchpad.net/~leonardr/beautifulsoup/bs4/view/head:/LICENSE', $null),
    @('F4', 'https://github.com/Ousret/charset_normalizer/blob/master/LICENSE', $null),
    @('F8', 'https://github.com/Mimino666/langdetect/blob/master/LICENSE', $null),
    @('F9', 'https://github.com/saffsd/langid.py/blob/master/LICENSE', $null),
    @('F11', 'https://github.com/lxml/lxml/blob/master/doc/licenses/BSD.txt', $null),
    @('F16', 'https://github.com/numpy/numpy/blob/master/LICENSE.txt', $null),
    @('F18', 'https://foss.heptapod.net/openpyxl/openpyxl/-/blob/branch/3.0/LICENCE.rst', $null),
    @('F26', 'https://github.com/python-openxml/python-docx/blob/master/LICENSE', $null),
    @('F28', 'https://github.com/requests/requests/blob/master/LICENSE', $null),
    @('F30', 'https://github.com/scipy/scipy/blob/master/LICENSE.txt', $null),
    @('F25', 'https://docs.python.org/3.8/license.html', 'psf-license-agreement-for-python-release'),
    @('F20', 'https://github.com/pyinstaller/pyinstaller/blob/develop/COPYING.txt', $null),
    @('F23', 'https://www.riverbankcomputing.com/static/Docs/PyQt5/introduction.html', 'license'),
    @('B25', 'https://www.python.org/', $null),
    @('B20', 'http://www.pyinstaller.org/', $null),
    @('B23', 'https://riverbankcomputing.com/software/pyqt/', $null),
    @('B12', 'https://matplotlib.org/', $null),
    @('B14', 'https://networkx.org/', $null),
    @('F12', 'https://matplotlib.org/users/license.html', $null),
    @('F14', 'https://github.com/networkx/networkx/blob/master/LICENSE.txt', $null),
    @('B38', 'https://github.com/amueller/word_cloud', $null),
    @('F38', 'https://github.com/amueller/word_cloud/blob/master/LICENSE', $null),
    @('B3', 'https://github.com/Esukhia/botok', $null),
    @('B5', 'https://github.com/cltk/cltk', $null),
    @('B7', 'https://github.com/fxsjy/jieba', $null),
    @('B13', 'https://github.com/taishi-i/nagisa', $null),
    @('B15', 'http://www.nltk.org/', $null),
    @('B17', 'https://github.com/yichen0831/opencc-python', $null),
    @('B19', 'https://github.com/lancopku/pkuseg-python', $null),
    @('B22', 'https://pyphen.org/', $null),
    @('B21', 'https://github.com/kmike/pymorphy2', $null),
    @('B24', 'https://github.com/PyThaiNLP/pythainlp', $null),
    @('B27', 'https://github.com/natasha/razdel', $null),
    @('B29', 'https://github.com/alvations/sacremoses', $null),
    @('B31', 'https://spacy.io/', $null),
    @('B32', 'https://github.com/ponrawee/ssg', $null),
    @('B34', 'https://github.com/fnl/syntok', $null),
    @('B35', 'https://github.com/sloria/TextBlob', $null),
    @('B36', 'https://github.com/mideind/Tokenizer', $null),
    @('B37', 'https://github.com/undertheseanlp/underthesea', $null),
    @('F3', 'https://github.com/Esukhia/botok/blob/master/LICENSE', $null),
    @('F5', 'https://github.com/cltk/cltk/blob/master/LICENSE', $null),
    @('F7', 'https://github.com/fxsjy/jieba/blob/master/LICENSE', $null),
    @('F13', 'https://github.com/taishi-i/nagisa/blob/master/LICENSE.txt', $null),
    @('F15', 'https://github.com/nltk/nltk/blob/develop/LICENSE.txt', $null),
    @('F17', 'https://github.com/yichen0831/opencc-python/blob/master/LICENSE.txt', $null),
    @('F19', 'https://github.com/lancopku/pkuseg-python/blob/master/LICENSE', $null),
    @('F22', 'https://github.com/Kozea/Pyphen/blob/master/LICENSE', $null),
    @('F21', 'https://github.com/kmike/pymorphy2/', 'pymorphy2'),
    @('F24', 'https://github.com/PyThaiNLP/pythainlp/blob/dev/LICENSE', $null),
    @('F27', 'https://github.com/natasha/razdel/blob/master/LICENSE', $null),
    @('F29', 'https://github.com/alvations/sacremoses/blob/master/LICENSE', $null),
    @('F31', 'https://github.com/explosion/spaCy/blob/master/LICENSE', $null),
    @('F32', 'https://github.com/ponrawee/ssg/blob/master/LICENSE', $null),
    @('F34', 'https://github.com/fnl/syntok/blob/master/LICENSE', $null),
    @('F35', 'https://github.com/sloria/TextBlob/blob/dev/LICENSE', $null),
    @('F36', 'https://github.com/mideind/Tokenizer/blob/master/LICENSE', $null),
    @('F37', 'https://github.com/undertheseanlp/underthesea/blob/master/LICENSE', $null),
    @('B6', 'https://github.com/Xangis/extra-stopwords', $null),
    @('B10', 'https://github.com/michmech/lemmatization-lists', $null),
    @('B33', 'https://github.com/stopwords-iso/stopwords-iso', $null),
    @('F6', 'https://github.com/Xangis/extra-stopwords/blob/master/LICENSE', $null),
    @('F10', 'https://github.com/michmech/lemmatization-lists/blob/master/LICENCE', $null),
    @('F33', 'https://github.com/stopwords-iso/stopwords-iso/blob/master/LICENSE', $null)

)

foreach ($item in $hyperlinkData) {
    $ref = $item[0]
    $address = $item[1]
    $subAddress = $item[2]
    $targetRange = $ws.Range($ref)
    if ($subAddress) {
        $ws.Hyperlinks.Add($targetRange, $address, $subAddress) | Out-Null
    } else {
        $ws.Hyperlinks.Add($targetRange, $address) | Out-Null
    }
}

# 4. Restore the selection / scroll position to reflect the new row layout
#    (mirrors where the author had scrolled to after removing the row).
$ws.Activate() | Out-Null
$ws.Range("D11").Select() | Out-Null
